$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 82
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 112
$ws.Range("E3").Value = 137
$ws.Range("H3").Value = 143
$ws.Range("J3").Value = 212
$ws.Range("K3").Value = 207
$ws.Range("L3").Value = 235
$ws.Range("E4").Value = 11
$ws.Range("J4").Value = 20
$ws.Range("D6").Value = 387
$ws.Range("E6").Value = 437
$ws.Range("F6").Value = 485
$ws.Range("G6").Value = 420
$ws.Range("I6").Value = 477
$ws.Range("J6").Value = 389
$ws.Range("K6").Value = 471
$ws.Range("L6").Value = 411
$ws.Range("D7").Value = 607
$ws.Range("E7").Value = 654
$ws.Range("F7").Value = 700
$ws.Range("G7").Value = 640
$ws.Range("H7").Value = 683
$ws.Range("I7").Value = 798
$ws.Range("J7").Value = 734
$ws.Range("K7").Value = 840
$ws.Range("L7").Value = 790

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("E2").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("E8").Value = 47
$ws.Range("F8").Value = 43
$ws.Range("G8").Value = 32
$ws.Range("J8").Value = 43
$ws.Range("L8").Value = 30
$ws.Range("F19").Value = 22
$ws.Range("L22").Value = 2
$ws.Range("L28").Value = 69
$ws.Range("G29").Value = 10
$ws.Range("K29").Value = 20
$ws.Range("E32").Value = 61
$ws.Range("K32").Value = 44
$ws.Range("E36").Value = 34
$ws.Range("K36").Value = 64
$ws.Range("E47").Value = 16
$ws.Range("E48").Value = 6
$ws.Range("D53").Value = 68
$ws.Range("E53").Value = 80
$ws.Range("H53").Value = 93
$ws.Range("I53").Value = 122
$ws.Range("L61").Value = 1
$ws.Range("J65").Value = 11
$ws.Range("J70").Value = 12
$ws.Range("J74").Value = 21
$ws.Range("L76").Value = 19
$ws.Range("D81").Value = 3
$ws.Range("E88").Value = 8
$ws.Range("D98").Value = 607
$ws.Range("E98").Value = 654
$ws.Range("F98").Value = 700
$ws.Range("G98").Value = 640
$ws.Range("H98").Value = 683
$ws.Range("I98").Value = 798
$ws.Range("J98").Value = 734
$ws.Range("K98").Value = 840
$ws.Range("L98").Value = 790

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 5
$ws.Range("L6").Value = 19

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("F2").Value = 8
$ws.Range("L3").Value = 13
$ws.Range("E6").Value = 37
$ws.Range("G6").Value = 23
$ws.Range("J6").Value = 23
$ws.Range("E7").Value = 47
$ws.Range("F7").Value = 43
$ws.Range("G7").Value = 32
$ws.Range("J7").Value = 43
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 16
$ws.Range("E6").Value = 48
$ws.Range("E7").Value = 61
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("E3").Value = 8
$ws.Range("K6").Value = 43
$ws.Range("E7").Value = 34
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 8

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H2").Value = 12
$ws.Range("I2").Value = 12
$ws.Range("H3").Value = 18
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 62
$ws.Range("I6").Value = 78
$ws.Range("D7").Value = 68
$ws.Range("E7").Value = 80
$ws.Range("H7").Value = 93
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 3

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 11

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 5
$ws.Range("J6").Value = 21

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("G5").Value = 8
$ws.Range("K5").Value = 15
$ws.Range("G6").Value = 10
$ws.Range("K6").Value = 20

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 22

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("E4").Value = 2
$ws.Range("E6").Value = 4

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 6

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("E5").Value = 10
$ws.Range("E6").Value = 16

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 12

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 2
